$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source cells are stored as literal text (inlineStr) in the workbook, e.g. "301.53",
# "-2.96%", "0.00006219". Force each target cell to Text format before assigning the
# new literal string so Excel does not auto-convert it to a number/percentage.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.96%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.52"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.10%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.073"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.77%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07994"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.75%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.942"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.39%"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.764"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.12%"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9266"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.05%"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1528"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "35.31%"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1905"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.57%"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08991"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.05%"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03454"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.31%"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09885"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.36%"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001397"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.50%"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005755"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.37%"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.532"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.85%"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.041"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.04%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.03%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3443"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.49%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.56%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.058"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.67%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2396"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.71%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04492"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.02%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001210"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.36%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004773"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.81%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-1.74%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-32.10%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01844"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.81%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04768"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.39%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01059"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "8.66%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007299"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.14%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1330"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.05%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.04%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009717"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-16.18%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006219"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.88%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.13%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.67"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-63.85%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.13%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.13%"
